$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.908.72'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.845.51'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.63%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.008'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.26'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4722'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3667'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +2.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07182'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.98%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +2.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.59'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.32%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.822.67'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.306'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.390'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.45'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.009'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008646'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.007'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.927.74'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.54'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.38%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.914'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.96%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.15'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.008'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '114.40'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.08%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.279'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7465'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.60%  '
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.165'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.93%  '
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.782'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.488'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.090'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.56%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +3.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01948'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.966'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.89%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5205'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.895'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.47%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.185'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.47'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4692'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.12%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.87'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.603'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '65.33'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06025'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8850'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.33%  '
